# Updated cryptos list - apply new Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.760.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.77"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7617"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.04"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3032"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.32"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06811"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07970"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.889.85"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7337"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.145"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.62"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.767.03"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.76"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.888"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.46"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007692"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.135.53"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.882"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.34"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.196"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1288"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.018"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.512"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.252"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05193"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.242"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7223"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.716"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01911"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.144"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4384"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.62"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.879"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8268"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.585"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.64"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.704"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.041.28"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.98"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.35%  "
